$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "'621"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = "'1443687.79"
$ws.Range("D2").Style = "Normal"

$ws.Range("C4").Value = "'1005"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = "'3535953.47"
$ws.Range("D4").Style = "Normal"

$ws.Range("C6").Value = "'636"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = "'2024112.78"
$ws.Range("D6").Style = "Normal"

$ws.Range("C8").Value = "'36"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "'154144.45"
$ws.Range("D8").Style = "Normal"

$ws.Range("C9").Value = "'168"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "'471180.44"
$ws.Range("D9").Style = "Normal"

$ws.Range("C14").Value = "'218"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'584362.00"
$ws.Range("D14").Style = "Normal"

$ws.Range("C16").Value = "'495"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'1823074.13"
$ws.Range("D16").Style = "Normal"

$ws.Range("C19").Value = "'8"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = "'28418.77"
$ws.Range("D19").Style = "Normal"

$ws.Range("C28").Value = "'276"
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = "'716237.45"
$ws.Range("D28").Style = "Normal"

$ws.Range("C35").Value = "'307"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'730119.71"
$ws.Range("D35").Style = "Normal"

$ws.Range("C36").Value = "'212"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'647703.04"
$ws.Range("D36").Style = "Normal"

$ws.Range("C37").Value = "'184"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'490860.14"
$ws.Range("D37").Style = "Normal"

$ws.Range("C39").Value = "'17"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'41492.85"
$ws.Range("D39").Style = "Normal"

$ws.Range("C45").Value = "'370"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'1007567.74"
$ws.Range("D45").Style = "Normal"

$ws.Range("C47").Value = "'609"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'2389882.99"
$ws.Range("D47").Style = "Normal"

$ws.Range("C48").Value = "'408"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'1380154.16"
$ws.Range("D48").Style = "Normal"

$ws.Range("C51").Value = "'3763"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'8676045.38"
$ws.Range("D51").Style = "Normal"

$ws.Range("C54").Value = "'4031"
$ws.Range("C54").Style = "Normal"
$ws.Range("D54").Value = "'14119639.38"
$ws.Range("D54").Style = "Normal"

$ws.Range("C74").Value = "'384"
$ws.Range("C74").Style = "Normal"
$ws.Range("D74").Value = "'968509.70"
$ws.Range("D74").Style = "Normal"

$ws.Range("C76").Value = "'919"
$ws.Range("C76").Style = "Normal"
$ws.Range("D76").Value = "'3192287.26"
$ws.Range("D76").Style = "Normal"

$ws.Range("C77").Value = "'519"
$ws.Range("C77").Style = "Normal"
$ws.Range("D77").Value = "'1705220.47"
$ws.Range("D77").Style = "Normal"

$ws.Range("C80").Value = "'393"
$ws.Range("C80").Style = "Normal"
$ws.Range("D80").Value = "'896342.96"
$ws.Range("D80").Style = "Normal"

$ws.Range("C92").Value = "'604"
$ws.Range("C92").Style = "Normal"
$ws.Range("D92").Value = "'1471849.94"
$ws.Range("D92").Style = "Normal"
